$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, shifting existing rows 233:254 down to 234:255.
$ws.Rows(233).Insert()

# Populate the newly inserted row 233 with the new weekly price record.
$ws.Range("A233").Value = 8
$ws.Range("B233").Value = "Terminal La Palmera de La Serena"
$ws.Range("C233").Value = "Coquimbo"
$ws.Range("D233").Value = 44769
$ws.Range("E233").Value = 4
$ws.Range("F233").Value = 100112031
$ws.Range("G233").Value = "Poroto verde"
$ws.Range("H233").Value = "Magnum"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 500
$ws.Range("K233").Value = 32000
$ws.Range("L233").Value = 33000
$ws.Range("M233").Value = 32500
$ws.Range("N233").Value = "$/malla 25 kilos"
$ws.Range("O233").Value = "Perú"
$ws.Range("P233").Value = 1300
$ws.Range("Q233").Value = 25
$ws.Range("R233").Value = "Hortaliza"
